$d = $word.ActiveDocument

# "Mostly" -> "Design: Mostly" also naturally covers "Mostly ok." -> "Design: Mostly ok."
# since Find/Replace matches the "Mostly" prefix of that string too.
$d.Content.Find.Execute("Mostly", $true, $false, $false, $false, $false, $true, 1, $false, "Design: Mostly", 2)

$d.Content.Find.Execute("Nothing special to mention", $true, $false, $false, $false, $false, $true, 1, $false, "Design: Nothing special to mention", 2)

$d.Content.Find.Execute("91xx Went well, some minor budget challenges", $true, $false, $false, $false, $false, $true, 1, $false, "Design: 91xx Went well, some minor budget challenges", 2)

$d.Content.Find.Execute("Internal communication ok. External communication with suppliers mostly ok.", $true, $false, $false, $false, $false, $true, 1, $false, "Design: Internal communication ok. External communication with suppliers mostly ok.", 2)

$d.Content.Find.Execute("Some things went to correct direction but regarding TK the opposite way.", $true, $false, $false, $false, $false, $true, 1, $false, "Design: Some things went to correct direction but regarding TK the opposite way.", 2)
